$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fill in the missing index number (A8) ---
$ws.Range("A8").Value = 7

# --- Row 9: now a fully populated acceptance-test row (was blank) ---
$ws.Rows.Item(9).RowHeight = 129.6
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Si la commande est validée par l'utilisateur, ce dernier est redirigé vers la page confirmation.html"
$ws.Range("C9").Value = "L'API renvoie en réponse, les informations entrées par l'utilisateur, les informations des articles commandés ainsi que le numéro de commande (orderId). L'orderId est passé dans l'URL de la page confirmation"
$ws.Range("D9").Value = "Dans la page confirmation.html, l'orderId est récupéré dans l'URL, et est affiché sur la page."
$ws.Range("E9").Value = "OK / Si la connexion avec l'API ne peut être établie, l'utilisateur en sera informé par un message d'alerte. "

# --- Row 10: another new acceptance-test row (was blank) ---
$ws.Rows.Item(10).RowHeight = 86.4
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Si la commande est passée, le localStorage est vidé de tous les articles,"
$ws.Range("C10").Value = "Le succès du passage de la commande provoque la redirection vers la page confirmation.html ainsi que la vidange du localStorage de tous les articles."
$ws.Range("D10").Value = "Reinitialisation du localStorage."
$ws.Range("E10").Value = "OK / Des problèmes de connexion avec l'API peuvent bloquer la commande."

# --- Row 11: another new acceptance-test row (was blank); A11 stays empty ---
$ws.Rows.Item(11).RowHeight = 64.8
$ws.Range("B11").Value = "Dans la page confirmation.html, suppression des informations concernant orderId"
$ws.Range("C11").Value = "Après récupération de orderId dans l'URL, orderId est effacé de l'instance de l'objet URL."
$ws.Range("D11").Value = "Suppression de orderId de l'instance de l'objet URL."
$ws.Range("E11").Value = "OK / "

# --- sheet view: zoom + scroll position + active selection ---
$win = $excel.ActiveWindow
$win.Zoom = 73
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("E11").Select()
